# ============================================================
# feat: add 2022-Q3 data
#
# Plan:
#  1. Workbook currently has sheets: 总计 (sheetId=1), 2022-Q2 (sheetId=2, holds
#     fund-detail rows for the single Q2 snapshot).
#  2. Target: 总计 (sheetId=1), 2022-Q3 (sheetId=2, NEW fund-detail rows),
#     2022-Q2 (sheetId=3, the ORIGINAL fund-detail rows, preserved verbatim).
#  3. 总计's own table gets a new leading row for 2022-Q3 and keeps the old
#     2022-Q2 row (shifted down).
# ============================================================

$wb = $excel.ActiveWorkbook
$zj = $wb.Worksheets.Item("总计")        # "总计" totals sheet
$source = $wb.Worksheets.Item("2022-Q2")  # current sheet holding the Q2 fund detail (to be cloned, then repurposed)

# ------------------------------------------------------------------
# Step 1: clone the current "2022-Q2" sheet (values+formats) into a brand new
# worksheet placed right after it -- that clone will remain named "2022-Q2"
# and keeps the exact original fund-detail data/styling untouched.
# ------------------------------------------------------------------
$newQ2 = $wb.Worksheets.Add($null, $source)
$source.Cells.Copy($newQ2.Range("A1"))
$newQ2.Range("A1").Clear()   # Copy() also stamped a blank A1 cell; the original sheet never had one
$newQ2.Name = "2022-Q2"

# ------------------------------------------------------------------
# Step 2: repurpose the original sheet (still physically "the 2022-Q2 sheet")
# as "2022-Q3" and overwrite its contents with the new Q3 fund-detail table.
# ------------------------------------------------------------------
$q3 = $source
$q3.Name = "2022-Q3"
$q3.Cells.Clear()

# Header row (style copied from the "总计" sheet's header cells so it matches the bold/border style)
$zj.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q3Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q3Headers.Length; $c++) {
    $q3.Cells.Item(1, 2 + $c).Value = $q3Headers[$c]
}

# Data rows 2..13. Column A holds a 0-based running index (style copied from
# the "总计" sheet's A2, which already carries that look). Columns D:G hold
# numeric-looking values that must stay literal text -> pre-format as Text.
$zj.Range("A2").Copy()
$q3.Range("A2:A13").PasteSpecial(-4122)

$q3.Range("B2:B13").NumberFormat = "@"
$q3.Range("D2:G13").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "008381"
$q3.Cells.Item(2, 3).Value = "前海开源新兴产业混合A"
$q3.Cells.Item(2, 4).Value = "7.73"
$q3.Cells.Item(2, 5).Value = "93.97"
$q3.Cells.Item(2, 6).Value = "4.73"
$q3.Cells.Item(2, 7).Value = "0.3656"
$q3.Cells.Item(2, 8).Value = 9

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "012442"
$q3.Cells.Item(3, 3).Value = "永赢稳健增长一年持有期混合E"
$q3.Cells.Item(3, 4).Value = "9.56"
$q3.Cells.Item(3, 5).Value = "26.04"
$q3.Cells.Item(3, 6).Value = "1.83"
$q3.Cells.Item(3, 7).Value = "0.1749"
$q3.Cells.Item(3, 8).Value = 3

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "011371"
$q3.Cells.Item(4, 3).Value = "华商远见价值混合型证券投资基金A"
$q3.Cells.Item(4, 4).Value = "3.85"
$q3.Cells.Item(4, 5).Value = "64.35"
$q3.Cells.Item(4, 6).Value = "4.00"
$q3.Cells.Item(4, 7).Value = "0.1540"
$q3.Cells.Item(4, 8).Value = 8

$q3.Cells.Item(5, 1).Value = 3
$q3.Cells.Item(5, 2).Value = "000800"
$q3.Cells.Item(5, 3).Value = "华商未来主题混合"
$q3.Cells.Item(5, 4).Value = "4.21"
$q3.Cells.Item(5, 5).Value = "72.31"
$q3.Cells.Item(5, 6).Value = "3.62"
$q3.Cells.Item(5, 7).Value = "0.1524"
$q3.Cells.Item(5, 8).Value = 7

$q3.Cells.Item(6, 1).Value = 4
$q3.Cells.Item(6, 2).Value = "001449"
$q3.Cells.Item(6, 3).Value = "华商双驱优选灵活配置混合"
$q3.Cells.Item(6, 4).Value = "2.41"
$q3.Cells.Item(6, 5).Value = "73.87"
$q3.Cells.Item(6, 6).Value = "4.11"
$q3.Cells.Item(6, 7).Value = "0.0991"
$q3.Cells.Item(6, 8).Value = 4

$q3.Cells.Item(7, 1).Value = 5
$q3.Cells.Item(7, 2).Value = "009932"
$q3.Cells.Item(7, 3).Value = "永赢稳健增长一年持有期混合A"
$q3.Cells.Item(7, 4).Value = "4.93"
$q3.Cells.Item(7, 5).Value = "26.04"
$q3.Cells.Item(7, 6).Value = "1.83"
$q3.Cells.Item(7, 7).Value = "0.0902"
$q3.Cells.Item(7, 8).Value = 3

$q3.Cells.Item(8, 1).Value = 6
$q3.Cells.Item(8, 2).Value = "008555"
$q3.Cells.Item(8, 3).Value = "华商龙头优势混合"
$q3.Cells.Item(8, 4).Value = "1.37"
$q3.Cells.Item(8, 5).Value = "77.58"
$q3.Cells.Item(8, 6).Value = "3.61"
$q3.Cells.Item(8, 7).Value = "0.0495"
$q3.Cells.Item(8, 8).Value = 8

$q3.Cells.Item(9, 1).Value = 7
$q3.Cells.Item(9, 2).Value = "002289"
$q3.Cells.Item(9, 3).Value = "华商改革创新股票A"
$q3.Cells.Item(9, 4).Value = "1.08"
$q3.Cells.Item(9, 5).Value = "79.89"
$q3.Cells.Item(9, 6).Value = "4.17"
$q3.Cells.Item(9, 7).Value = "0.0450"
$q3.Cells.Item(9, 8).Value = 5

$q3.Cells.Item(10, 1).Value = 8
$q3.Cells.Item(10, 2).Value = "014729"
$q3.Cells.Item(10, 3).Value = "前海开源新兴产业混合C"
$q3.Cells.Item(10, 4).Value = "0.60"
$q3.Cells.Item(10, 5).Value = "93.97"
$q3.Cells.Item(10, 6).Value = "4.73"
$q3.Cells.Item(10, 7).Value = "0.0284"
$q3.Cells.Item(10, 8).Value = 9

$q3.Cells.Item(11, 1).Value = 9
$q3.Cells.Item(11, 2).Value = "010403"
$q3.Cells.Item(11, 3).Value = "华商景气优选混合"
$q3.Cells.Item(11, 4).Value = "0.61"
$q3.Cells.Item(11, 5).Value = "76.85"
$q3.Cells.Item(11, 6).Value = "4.19"
$q3.Cells.Item(11, 7).Value = "0.0256"
$q3.Cells.Item(11, 8).Value = 6

$q3.Cells.Item(12, 1).Value = 10
$q3.Cells.Item(12, 2).Value = "016052"
$q3.Cells.Item(12, 3).Value = "华商改革创新股票C"
$q3.Cells.Item(12, 4).Value = "0.32"
$q3.Cells.Item(12, 5).Value = "79.89"
$q3.Cells.Item(12, 6).Value = "4.17"
$q3.Cells.Item(12, 7).Value = "0.0133"
$q3.Cells.Item(12, 8).Value = 5

$q3.Cells.Item(13, 1).Value = 11
$q3.Cells.Item(13, 2).Value = "011372"
$q3.Cells.Item(13, 3).Value = "华商远见价值混合型证券投资基金C"
$q3.Cells.Item(13, 4).Value = "0.28"
$q3.Cells.Item(13, 5).Value = "64.35"
$q3.Cells.Item(13, 6).Value = "4.00"
$q3.Cells.Item(13, 7).Value = "0.0112"
$q3.Cells.Item(13, 8).Value = 8

# Match the page margins used by the rest of the workbook's primary sheets
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# Step 3: update the "总计" (totals) sheet -- insert the 2022-Q3 total as the new
# first data row, and keep the 2022-Q2 total as the following row.
# ------------------------------------------------------------------
$zj.Range("A2").Copy()
$zj.Range("A3").PasteSpecial(-4122)   # xlPasteFormats, so A3 gets the same index style as A2

$zj.Range("B3").Value = $zj.Range("B2").Value()
$zj.Range("C3").Value = $zj.Range("C2").Value()
$zj.Range("D3").Value = $zj.Range("D2").Value()
$zj.Range("A3").Value = 1

$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 12
$zj.Range("D2").Value = 1.21

Write-Output "2022-Q3 sheet added"
